$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D7").Value = "implemented"
$ws.Range("E7").Value = "Execution backtest implemented: base portfolio run selection + ideal vs realistic fills (CLOSE vs NEXT_OPEN) with slippage/charges + comparison chart and API tests."
$ws.Range("F7").Value = "27/12/2025 04:42"
